# Updated: st 12. 05. 2021
# Applies corrections to AgTests (F) and AgPosit (G) columns for rows 393-432
# matching the diff of OpenData_Slovakia_Covid_DailyStats.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Cell = "F393"; Value = 307757 },
    @{ Cell = "G393"; Value = 1239 },
    @{ Cell = "F395"; Value = 750944 },
    @{ Cell = "F401"; Value = 273382 },
    @{ Cell = "G401"; Value = 936 },
    @{ Cell = "F402"; Value = 717944 },
    @{ Cell = "G403"; Value = 733 },
    @{ Cell = "F404"; Value = 224784 },
    @{ Cell = "G404"; Value = 910 },
    @{ Cell = "F405"; Value = 174319 },
    @{ Cell = "G405"; Value = 695 },
    @{ Cell = "F406"; Value = 170456 },
    @{ Cell = "G406"; Value = 681 },
    @{ Cell = "F408"; Value = 303391 },
    @{ Cell = "F410"; Value = 363024 },
    @{ Cell = "G410"; Value = 629 },
    @{ Cell = "F411"; Value = 224998 },
    @{ Cell = "F412"; Value = 175814 },
    @{ Cell = "F413"; Value = 149021 },
    @{ Cell = "G413"; Value = 659 },
    @{ Cell = "F414"; Value = 146995 },
    @{ Cell = "F415"; Value = 304572 },
    @{ Cell = "F416"; Value = 660105 },
    @{ Cell = "F417"; Value = 332533 },
    @{ Cell = "F418"; Value = 200583 },
    @{ Cell = "G418"; Value = 700 },
    @{ Cell = "F419"; Value = 147468 },
    @{ Cell = "F420"; Value = 136698 },
    @{ Cell = "G420"; Value = 491 },
    @{ Cell = "F421"; Value = 150455 },
    @{ Cell = "F422"; Value = 293849 },
    @{ Cell = "G422"; Value = 641 },
    @{ Cell = "F425"; Value = 136506 },
    @{ Cell = "F426"; Value = 105024 },
    @{ Cell = "F428"; Value = 100116 },
    @{ Cell = "F429"; Value = 171464 },
    @{ Cell = "G429"; Value = 446 },
    @{ Cell = "F430"; Value = 168503 },
    @{ Cell = "G430"; Value = 271 },
    @{ Cell = "F431"; Value = 162725 },
    @{ Cell = "G431"; Value = 391 },
    @{ Cell = "F432"; Value = 120256 },
    @{ Cell = "G432"; Value = 419 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
